# Add new "Save" column (H) to the s_vals sheet, matching style of existing
# header cells and plain numeric style of existing data cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# H1: header cell "Save" -- copy formatting from the adjacent header cell (G1)
# then overwrite the value/text.
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# H2: data cell with numeric value 0 (default, unstyled like the other data cells)
$ws.Range("H2").Value = 0
